$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated fitting parameters
$ws.Range("J2").Value = 0.01409
$ws.Range("K2").Value = 0.11079

# Mark the K3 "-" unit cell (ready to run / detail tests) with an underline
$ws.Range("K3").Font.Underline = $true

# Move the selection to K3
$ws.Activate()
$ws.Range("K3").Select()
